# "Ajout de la fonction favorite" -- Sprint Burndown Chart updates:
# points re-estimated and status/owner filled in for several tasks
# related to the "favorite route" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Burndown Chart")
$ws.Activate()

# Row 9 - "Buttons annuler et confirmer": now Closed, owned by Fabio
$ws.Range("F9").Value2 = "Closed"
$ws.Range("G9").Value2 = "Fabio"

# Row 10 - "Rédirection page login.php": now Closed, owned by Fabio
$ws.Range("F10").Value2 = "Closed"
$ws.Range("G10").Value2 = "Fabio"

# Row 12 - Us_005 "Ajouter le parcours aux favoris": points 5 -> 10, now Closed
$ws.Range("C12").Value2 = 10
$ws.Range("F12").Value2 = "Closed"

# Row 13 - Us_006 "Consulter ses informations": points 5 -> 15
$ws.Range("C13").Value2 = 15

# Row 14 - "Modifier les informations": points 5 -> 15
$ws.Range("C14").Value2 = 15

# Row 15 - Us_007 "Consulter le tableau des favoris": points 5 -> 15, now In progress, owned by Alex
$ws.Range("C15").Value2 = 15
$ws.Range("F15").Value2 = "In progress"
$ws.Range("G15").Value2 = "Alex"

# Row 16 - "Supprimer un des parcours favoris": points 5 -> 10, owned by Alex
$ws.Range("C16").Value2 = 10
$ws.Range("G16").Value2 = "Alex"

# Row 17 - Us_008 "Supprimer un compte": points 5 -> 10, now In progress, owned by Damiano
$ws.Range("C17").Value2 = 10
$ws.Range("F17").Value2 = "In progress"
$ws.Range("G17").Value2 = "Damiano"

# Row 18 - Us_009 "Ajout d'un parcours": points 5 -> 25, now In progress, owned by Damiano
$ws.Range("C18").Value2 = 25
$ws.Range("F18").Value2 = "In progress"
$ws.Range("G18").Value2 = "Damiano"

# Match the author's final view state: zoomed in, looking at F13
$excel.ActiveWindow.Zoom = 145
$ws.Range("F13").Select()
